# Insert a new data row before row 147 (pushing existing rows 147-171 down
# to 148-172) and populate the new row 147 with the latest price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(147).Insert()

$ws.Cells.Item(147, 1).Value = 10
$ws.Cells.Item(147, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(147, 3).Value = "La Araucanía"
$ws.Cells.Item(147, 4).Value = 44951
$ws.Cells.Item(147, 5).Value = 9
$ws.Cells.Item(147, 6).Value = 100112031
$ws.Cells.Item(147, 7).Value = "Poroto verde"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 40
$ws.Cells.Item(147, 11).Value = 30000
$ws.Cells.Item(147, 12).Value = 30000
$ws.Cells.Item(147, 13).Value = 30000
$ws.Cells.Item(147, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(147, 15).Value = "Región del Maule"
$ws.Cells.Item(147, 16).Value = 1200
$ws.Cells.Item(147, 17).Value = 25
$ws.Cells.Item(147, 18).Value = "Hortaliza"
